$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 437, shifting the existing rows
# (old 437-458) down to (439-460). This mirrors the source data picking
# up a new weekly price entry (dated 45041) for "Apio" at
# "Vega Monumental Concepción", while the oldest entries that fall off
# the bottom of the 24-row rolling window (previously rows 457-458)
# reappear duplicated at the new end of the range (rows 459-460).
$ws.Rows.Item(437).Insert()
$ws.Rows.Item(437).Insert()

$ws.Range("A437:A438").Value = 11
$ws.Range("B437:B438").Value = "Vega Monumental Concepción"
$ws.Range("C437:C438").Value = "Bíobío"
$ws.Range("D437:D438").Value = 45041
$ws.Range("E437:E438").Value = 8
$ws.Range("F437:F438").Value = 100112017
$ws.Range("G437:G438").Value = "Apio"
$ws.Range("H437:H438").Value = "Americana (o)"

$ws.Range("I437").Value = "Primera"
$ws.Range("J437").Value = 270
$ws.Range("K437").Value = 8500
$ws.Range("L437").Value = 9000
$ws.Range("M437").Value = 8778
$ws.Range("P437").Value = 1463

$ws.Range("I438").Value = "Segunda"
$ws.Range("J438").Value = 220
$ws.Range("K438").Value = 7000
$ws.Range("L438").Value = 7500
$ws.Range("M438").Value = 7273
$ws.Range("P438").Value = 1212

$ws.Range("N437:N438").Value = "$/docena de matas"
$ws.Range("O437:O438").Value = "Región de Coquimbo"
$ws.Range("Q437:Q438").Value = 6
$ws.Range("R437:R438").Value = "Hortaliza"
